$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (changed) date column C for rows 2 through 14
# from serial date 45204 (2023-10-05) to 45207 (2023-10-08),
# keeping the existing date number formatting on the cells.
for ($row = 2; $row -le 14; $row++) {
    $ws.Cells.Item($row, 3).Value = 45207
}
